# Tarefa de estatística 01 - Tabela de Frequencias
# Add todos os cálculos de mediana, média e moda
#
# - Column D (fp(%)) and column F (fac(%)) currently hold literal text
#   strings such as "4.0%" / "28.000000000000004%". Replace them with the
#   real numeric fraction (fp, and fac/total respectively) formatted with
#   a "##%" custom number format, so they are proper percentage values.
# - Column I (moda) gets the data values that occur most often
#   (frequency = 2): 30, 44 and 87, in rows 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$total = 25

# row -> (fp value for column D, cumulative fac count for column F)
$rows = @(
    @{ r = 3;  fp = 0.04; fac = 1 },
    @{ r = 4;  fp = 0.04; fac = 2 },
    @{ r = 5;  fp = 0.04; fac = 3 },
    @{ r = 6;  fp = 0.04; fac = 4 },
    @{ r = 7;  fp = 0.04; fac = 5 },
    @{ r = 8;  fp = 0.08; fac = 7 },
    @{ r = 9;  fp = 0.04; fac = 8 },
    @{ r = 10; fp = 0.08; fac = 10 },
    @{ r = 11; fp = 0.04; fac = 11 },
    @{ r = 12; fp = 0.04; fac = 12 },
    @{ r = 13; fp = 0.04; fac = 13 },
    @{ r = 14; fp = 0.04; fac = 14 },
    @{ r = 15; fp = 0.04; fac = 15 },
    @{ r = 16; fp = 0.04; fac = 16 },
    @{ r = 17; fp = 0.04; fac = 17 },
    @{ r = 18; fp = 0.04; fac = 18 },
    @{ r = 19; fp = 0.04; fac = 19 },
    @{ r = 20; fp = 0.08; fac = 21 },
    @{ r = 21; fp = 0.04; fac = 22 },
    @{ r = 22; fp = 0.04; fac = 23 },
    @{ r = 23; fp = 0.04; fac = 24 },
    @{ r = 24; fp = 0.04; fac = 25 }
)

foreach ($row in $rows) {
    $r = $row.r

    # Column D: fp(%) -- same underlying fraction as column C (fp), now
    # stored as a real number displayed with the "##%" number format
    # instead of a literal text string.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $row.fp
    $dCell.NumberFormat = "##%"

    # Column F: fac(%) -- cumulative frequency (fac) divided by the total
    # number of observations, displayed with the "##%" number format.
    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Value = $row.fac / $total
    $fCell.NumberFormat = "##%"
}

# Column I (rows 3-5): "moda" -- the data values that occur most often
# (frequency = 2): 30, 44 and 87.
$ws.Cells.Item(3, 9).Value = 30
$ws.Cells.Item(4, 9).Value = 44
$ws.Cells.Item(5, 9).Value = 87
